# 🚌 141: 31/12 16:51 LP1912+6203+6173
# Update scrape timestamp / row-count headers and append newly-scraped rows
# to the "LP1912", "LP1912-215" and "6203-6173" sheets.

$wb = $excel.ActiveWorkbook

$newStamp = "Última actualización: 31/12/2025 13:51:30"

# ---------------------------------------------------------------------------
# Sheet "LP1912": columns A(title/blank) B=Hora_Scrap C=Hora_Llegada
#                 D=Linea E=Minutos F=Parada G=Fecha
# ---------------------------------------------------------------------------
$wsLP = $wb.Worksheets.Item("LP1912")

$wsLP.Cells.Item(2, 1).Value = $newStamp
$wsLP.Cells.Item(3, 1).Value = "Total filas: 989"

$lpRows = @(
    @("13:51:19", "14:01", "16_SANTA ANA",        10, "LP1912", "31/12/2025"),
    @("13:51:19", "14:01", "17_ROMERO",            10, "LP1912", "31/12/2025"),
    @("13:51:19", "14:03", "23_HERNANDEZ",         12, "LP1912", "31/12/2025"),
    @("13:51:19", "14:13", "16_SANTA ANA",         22, "LP1912", "31/12/2025"),
    @("13:51:19", "14:25", "11_ETCHEVERRY",        34, "LP1912", "31/12/2025"),
    @("13:51:19", "14:25", "16_SANTA ANA",         34, "LP1912", "31/12/2025"),
    @("13:51:19", "14:37", "16_P MOR-SANTA ANA",   46, "LP1912", "31/12/2025"),
    @("13:51:19", "14:40", "17X38_ROMERO",         49, "LP1912", "31/12/2025"),
    @("13:51:19", "14:41", "23_HERNANDEZ",         50, "LP1912", "31/12/2025"),
    @("13:51:19", "15:13", "15_ABASTO",            82, "LP1912", "31/12/2025"),
    @("13:51:19", "15:14", "10_OLMOS",             83, "LP1912", "31/12/2025"),
    @("13:51:19", "15:25", "11_ETCHEVERRY",        94, "LP1912", "31/12/2025")
)

$startRow = 979
for ($i = 0; $i -lt $lpRows.Count; $i++) {
    $r = $startRow + $i
    $row = $lpRows[$i]
    $wsLP.Cells.Item($r, 2).Value = $row[0]
    $wsLP.Cells.Item($r, 3).Value = $row[1]
    $wsLP.Cells.Item($r, 4).Value = $row[2]
    $wsLP.Cells.Item($r, 5).Value = $row[3]
    $wsLP.Cells.Item($r, 6).Value = $row[4]
    $wsLP.Cells.Item($r, 7).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": only the "last updated" stamp changes
# ---------------------------------------------------------------------------
$wsLP215 = $wb.Worksheets.Item("LP1912-215")
$wsLP215.Cells.Item(2, 1).Value = $newStamp

# ---------------------------------------------------------------------------
# Sheet "6203-6173": columns A(title/blank) B=Fecha C=Hora_Scrap
#                     D=Hora_Llegada E=Linea F=Minutos G=Parada
# ---------------------------------------------------------------------------
$ws6203 = $wb.Worksheets.Item("6203-6173")

$ws6203.Cells.Item(2, 1).Value = $newStamp
$ws6203.Cells.Item(3, 1).Value = "Total filas: 121"

$c6203Rows = @(
    @("31/12/2025", "13:51:25", "14:34", "215C_LA PLATA", 43, "L6203"),
    @("31/12/2025", "13:51:30", "15:00", "215A_LA PLATA", 69, "L6173")
)

$startRow6203 = 121
for ($i = 0; $i -lt $c6203Rows.Count; $i++) {
    $r = $startRow6203 + $i
    $row = $c6203Rows[$i]
    $ws6203.Cells.Item($r, 2).Value = $row[0]
    $ws6203.Cells.Item($r, 3).Value = $row[1]
    $ws6203.Cells.Item($r, 4).Value = $row[2]
    $ws6203.Cells.Item($r, 5).Value = $row[3]
    $ws6203.Cells.Item($r, 6).Value = $row[4]
    $ws6203.Cells.Item($r, 7).Value = $row[5]
}
